$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new value, as captured from the upstream coinranking refresh
$updates = [ordered]@{
    'D2' = '52.275.06'
    'E2' = '  +1.20%  '
    'D3' = '2.823.45'
    'E3' = '  +2.85%  '
    'E4' = '  -0.06%  '
    'D5' = '355.03'
    'E5' = '  +6.58%  '
    'D6' = '114.19'
    'E6' = '  -1.95%  '
    'E7' = '  +2.50%  '
    'E8' = '  -0.01%  '
    'D9' = '0.604'
    'E9' = '  +5.88%  '
    'D10' = '42.01'
    'E10' = '  +1.40%  '
    'D11' = '0.0851'
    'E11' = '  -1.22%  '
    'D12' = '20.14'
    'E12' = '  +0.16%  '
    'E13' = '  +1.32%  '
    'D14' = '7.78'
    'E14' = '  +2.99%  '
    'D15' = '3.253.75'
    'E15' = '  +2.47%  '
    'D16' = '2.826.98'
    'E16' = '  +3.12%  '
    'D17' = '0.891'
    'E17' = '  +1.40%  '
    'D18' = '52.160.66'
    'E18' = '  +1.12%  '
    'D19' = '3.18'
    'E19' = '  +1.65%  '
    'B20' = 'Uniswap'
    'C20' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D20' = '7.30'
    'E20' = '  +6.99%  '
    'B21' = 'InternetComputer(DFINITY)'
    'C21' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D21' = '13.81'
    'E21' = '  +2.50%  '
    'D22' = '0.0₃0996'
    'E22' = '  +1.91%  '
    'D23' = '270.83'
    'E23' = '  -3.09%  '
    'D24' = '69.65'
    'E24' = '  +0.27%  '
    'E25' = '  +5.38%  '
    'D26' = '26.72'
    'E26' = '  +0.12%  '
    'E27' = '  +0.01%  '
    'D28' = '10.27'
    'E28' = '  +1.04%  '
    'E29' = '  +1.40%  '
    'E30' = '  +0.70%  '
    'D31' = '50.75'
    'E31' = '  +1.52%  '
    'D32' = '33.90'
    'E32' = '  -3.13%  '
    'D33' = '5.87'
    'E33' = '  +6.02%  '
    'D34' = '0.0445'
    'E34' = '  +28.86%  '
    'E35' = '  +1.89%  '
    'E36' = '  -0.14%  '
    'E37' = '  +1.28%  '
    'D38' = '4.90'
    'E38' = '  -1.46%  '
    'D39' = '18.43'
    'E39' = '  -2.69%  '
    'E40' = '  +1.39%  '
    'D41' = '2.60'
    'E41' = '  +9.57%  '
    'D42' = '128.31'
    'E42' = '  +0.16%  '
    'D43' = '23.52'
    'E43' = '  +2.43%  '
    'E44' = '  +1.98%  '
    'E45' = '  +1.49%  '
    'E46' = '  +1.77%  '
    'D47' = '2.042.06'
    'E47' = '  -2.42%  '
    'E48' = '  +3.23%  '
    'D49' = '0.971'
    'E49' = '  +12.11%  '
    'E50' = '  +3.21%  '
    'D51' = '60.37'
    'E51' = '  +1.00%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text interpretation so numeric-looking strings (e.g. "42.01", "1.00")
    # and formatted values (e.g. "  +1.20%  ") are preserved exactly as-is.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}

